$d = $word.ActiveDocument

# --- Paragraph 6 (Objetivos content) ---
$p = $d.Paragraphs.Item(6).Range
$null = $p.Find.Execute("Oferecer ao aluno uma visão geral dos conceitos fundamentais e das fases do processo de desenvolvimento e Implementação de Sistemas de Informação no sentido de capacitá-lo analisa e projetar tais sistemas", $true, $false, $false, $false, $false, $true, 1, $false, "Sistemas de Informação. Projeto de Sistemas de Informação. Tecnologia de Informação. Processo de Desenvolvimento de Sistema de Informação.", 2)

# --- Paragraph 8 (Docente ListBullet content) ---
$p = $d.Paragraphs.Item(8).Range
$null = $p.Find.Execute("5840917 - Fabricio Maciel Gomes", $true, $false, $false, $false, $false, $true, 1, $false, "Oferecer ao aluno uma visão geral dos conceitos fundamentais e das fases do processo de desenvolvimento e Implementação de Sistemas de Informação no sentido de capacitá-lo analisa e projetar tais sistemas", 2)

# --- Paragraph 10 (content under "Programa resumido") ---
$p = $d.Paragraphs.Item(10).Range
$null = $p.Find.Execute("Sistemas de Informação. Projeto de Sistemas de Informação. Tecnologia de Informação. Processo de Desenvolvimento de Sistema de Informação.", $true, $false, $false, $false, $false, $true, 1, $false, "1. Sistemas de Informação`v1.1. Sistemas de Processamento de Informações;`v1.2. Sistemas de Informações Gerenciais;`v1.3. Sistema de Apoio à Decisão;`v1.4. Sistemas de Informação no Comércio Eletrônico;`v1.5. Sistemas de Informação em Cadeia de Suprimentos;`v1.6. Sistemas inteligentes nos negócios;`v1.7. Sistemas estratégicos. `v2. Projeto de Sistemas de Informação.`v2.1. Especificação das Saídas;`v2.2. Especificação dos Arquivos;`v2.3. Especificação das Entradas;`v2.4. Especificação do Processamento.`v3. Tecnologia de Informação.`v3.1. Evolução da Computação;`v3.2. Recursos Computacionais.`v4. Processo de Desenvolvimento de Sistemas de Informação.`v4.1. Definição do Negócio;`v4.2. Identificação do Problema e/ou Oportunidades;`v4.3. Seleção do Sistema de Informação;`v4.4. Implementação do Sistema de Informação;`v4.5. Avaliação da Eficácia do Sistema de Informação;", 2)

# --- Paragraph 12 (content under "Programa") ---
$p = $d.Paragraphs.Item(12).Range
$null = $p.Find.Execute("1. Sistemas de Informação`v1.1. Sistemas de Processamento de Informações;`v1.2. Sistemas de Informações Gerenciais;`v1.3. Sistema de Apoio à Decisão;`v1.4. Sistemas de Informação no Comércio Eletrônico;`v1.5. Sistemas de Informação em Cadeia de Suprimentos;`v1.6. Sistemas inteligentes nos negócios;`v1.7. Sistemas estratégicos. `v2. Projeto de Sistemas de Informação.`v2.1. Especificação das Saídas;`v2.2. Especificação dos Arquivos;`v2.3. Especificação das Entradas;`v2.4. Especificação do Processamento.`v3. Tecnologia de Informação.`v3.1. Evolução da Computação;`v3.2. Recursos Computacionais.`v4. Processo de Desenvolvimento de Sistemas de Informação.`v4.1. Definição do Negócio;`v4.2. Identificação do Problema e/ou Oportunidades;`v4.3. Seleção do Sistema de Informação;`v4.4. Implementação do Sistema de Informação;`v4.5. Avaliação da Eficácia do Sistema de Informação;", $true, $false, $false, $false, $false, $true, 1, $false, "Aulas expositivas teóricas, aulas práticas, aulas de exercícios.", 2)

# --- Paragraph 14 (Avaliacao block): replace in reverse order to avoid collisions ---
$p = $d.Paragraphs.Item(14).Range
$null = $p.Find.Execute("Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação.", $true, $false, $false, $false, $false, $true, 1, $false, "HAL R. VARIAN, H. R.; FARRELL, J., SHAPIRO, C. The economics of information technology: an introduction. Cambridge University Press, 2004`vLAURINDO, F. J. B.; ROTONDARO, R. G. orgs. Gestão integrada de processos e da tecnologia da informação. São Paulo:Atlas, 2006.`vLAURINDO, F.J.B. Tecnologia da Informação: Eficácia nas Organizações. São Paulo, Editora Futura, 2002.`vSTAIR, R.M., Princípios de Sistema de Informação: Uma Abordagem Gerencial, Rio de Janeiro, LTC, 1998.`vTURBAN, E. et al. Information Technology for Management: Transforming Organizations in the Digital Economy. 7th edition, Wiley, 2009.`vTURBAN, E., RAIANER JR, K., POTTER, R. E., Administração de Tecnologia da Informação: Teoria e Prática”, São Paulo, Editora Campus, 2003.", 2)
$p = $d.Paragraphs.Item(14).Range
$null = $p.Find.Execute("Média Aritmética das atividades avaliativas realizadas.`v", $true, $false, $false, $false, $false, $true, 1, $false, "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação.`v", 2)
$p = $d.Paragraphs.Item(14).Range
$null = $p.Find.Execute("Aulas expositivas teóricas, aulas práticas, aulas de exercícios.`v", $true, $false, $false, $false, $false, $true, 1, $false, "Média Aritmética das atividades avaliativas realizadas.`v", 2)

# --- Paragraph 16 (Bibliografia content) ---
$p = $d.Paragraphs.Item(16).Range
$null = $p.Find.Execute("HAL R. VARIAN, H. R.; FARRELL, J., SHAPIRO, C. The economics of information technology: an introduction. Cambridge University Press, 2004`vLAURINDO, F. J. B.; ROTONDARO, R. G. orgs. Gestão integrada de processos e da tecnologia da informação. São Paulo:Atlas, 2006.`vLAURINDO, F.J.B. Tecnologia da Informação: Eficácia nas Organizações. São Paulo, Editora Futura, 2002.`vSTAIR, R.M., Princípios de Sistema de Informação: Uma Abordagem Gerencial, Rio de Janeiro, LTC, 1998.`vTURBAN, E. et al. Information Technology for Management: Transforming Organizations in the Digital Economy. 7th edition, Wiley, 2009.`vTURBAN, E., RAIANER JR, K., POTTER, R. E., Administração de Tecnologia da Informação: Teoria e Prática”, São Paulo, Editora Campus, 2003.", $true, $false, $false, $false, $false, $true, 1, $false, "5840917 - Fabricio Maciel Gomes", 2)

Write-Output "done"
